$d = $word.ActiveDocument
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# 1) Split "Debut de la realisation du CSS sur le front office." into three
#    runs ("Debut de la " / "reflexion et la " / "realisation ...") and drop
#    the _GoBack bookmark that currently sits at the end of that paragraph
#    (it gets re-created further down, attached to a different paragraph).
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

$rng = $d.Content
$rng.Find.Execute("Début de la réalisation du CSS sur le front office.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = ""
$rng.InsertAfter("Début de la ")
$rng.Collapse(0)
$rng.InsertAfter("réflexion et la ")
$rng.Collapse(0)
$rng.InsertAfter("réalisation du CSS sur le front office.")
$rng.Collapse(0)

# ---------------------------------------------------------------------------
# 2) The trailing "Documentation sur les noms de session (hors travail)."
#    paragraph (previously split across two runs "D" / "ocumentation...")
#    is replaced by itself (merged into one run) plus a whole new block of
#    journal content that follows it, up to the end of the document.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$target = $lastPara.Range

$xml = @"
<w:p $ns>
  <w:pPr><w:spacing w:after="0"/></w:pPr>
  <w:r><w:t>Documentation sur les noms de session (hors travail).</w:t></w:r>
</w:p>
<w:p $ns>
  <w:pPr><w:spacing w:after="0"/></w:pPr>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr><w:b/><w:u w:val="single"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:b/><w:u w:val="single"/></w:rPr>
    <w:t>03/05/2018</w:t>
  </w:r>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr><w:b/><w:u w:val="single"/></w:rPr>
  </w:pPr>
</w:p>
<w:p $ns>
  <w:pPr><w:spacing w:after="0"/></w:pPr>
  <w:r><w:t>Avancement dans le sprint 2, avancement dans la réflexion et la réalisation du CSS.</w:t></w:r>
</w:p>
<w:p $ns>
  <w:pPr><w:spacing w:after="0"/></w:pPr>
</w:p>
<w:p $ns>
  <w:pPr><w:spacing w:after="0"/></w:pPr>
  <w:r><w:t xml:space="preserve">Entretien avec Françoise PHILIBERT, </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:b/></w:rPr><w:t>product</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:b/></w:rPr><w:t>owner</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t>et membre du service communication de l’entreprise : mise au point pour évaluer les points positifs et négatifs des fonctionnalités implémentées.</w:t></w:r>
</w:p>
<w:p $ns>
  <w:pPr><w:spacing w:after="0"/></w:pPr>
</w:p>
<w:p $ns>
  <w:pPr><w:spacing w:after="0"/></w:pPr>
  <w:r><w:t xml:space="preserve">Bilan : bon dans l’ensemble, quelques micro-correctifs à ajouter (augmenter la taille des </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>textareas</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> par exemple, modification des messages envoyés par mail) et annonce de nouvelles fonctionnalités :</w:t></w:r>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:pStyle w:val="Paragraphedeliste"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
    <w:spacing w:after="0"/>
  </w:pPr>
  <w:r><w:t>Possibilité de trier par nom de structure et par date les différents projets</w:t></w:r>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:pStyle w:val="Paragraphedeliste"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
    <w:spacing w:after="0"/>
  </w:pPr>
  <w:r><w:t>Possibilité les projets triés, ou non, ou les deux.</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:pStyle w:val="Paragraphedeliste"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
    <w:spacing w:after="0"/>
  </w:pPr>
  <w:r><w:t>Ajout d’un onglet de suivi contenant des renseignements spécifiques permettant un meilleur suivi des projets (contenu divulgué dans le futur)</w:t></w:r>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:pStyle w:val="Paragraphedeliste"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
    <w:spacing w:after="0"/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Possibilité de </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>upload</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> et download des fichiers personnels.</w:t></w:r>
</w:p>
"@

$target.InsertXML($xml)
